$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 179
$ws.Range("B1").Value = 245.6000000000004
$ws.Range("C1").Value = 284.8000000000011

$ws.Range("A2").Value = 179
$ws.Range("B2").Value = 257.4000000000015
$ws.Range("C2").Value = 179
